$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 10: the long "Objetivos" paragraph is replaced by the docente name
# ------------------------------------------------------------------
$ws.Range("B10").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C10").Value = "4808662 - Lucrécio Fábio dos Santos"

# ------------------------------------------------------------------
# Row 13 - label "Programa resumido:" gains an A cell, B/C become "Semestral"
# ------------------------------------------------------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# ------------------------------------------------------------------
# Row 14 - label shifts to "Short syllabus:", values become the short syllabus text
# ------------------------------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1) Determination of Reynolds number 2) Liquids flow measurements 3) Determination of the fluids velocity profile in pipe 4) Liquid flow from a cylindrical tank 5) Determination of friction factor in pipes 6) Rheological tests."
$ws.Range("C14").Value = "1) Determination of Reynolds number 2) Liquids flow measurements 3) Determination of the fluids velocity profile in pipe 4) Liquid flow from a cylindrical tank 5) Determination of friction factor in pipes 6) Rheological tests."

# ------------------------------------------------------------------
# Row 15 - label becomes "Programa:", values become "01/01/2022" (must stay TEXT,
# not be auto-converted to a date serial number) -- copy the already-text value
# that lives in B8/C8 so type + number format survive untouched.
# ------------------------------------------------------------------
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Row 16 - label becomes "Syllabus:", values become the full syllabus text
# ------------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1) Determination of the Reynolds number using the experimental apparatus. 2) Liquids flow measurements: Determining the flow by using a rotameter, orifice plate and Venturi. 3) Liquid flow from a cylindrical tank: verification of a mathematical model based on the conservation equations of mass and energy to determine the time of emptying reservoir and compared the results with experimental data. 4) Determination of friction factor in pipes: evaluation of the friction factor as a function of Reynolds number in pipes. 5) Rheological tests: with the use of various types of viscometers to determine the dynamic and kinematic viscosities."
$ws.Range("C16").Value = "1) Determination of the Reynolds number using the experimental apparatus. 2) Liquids flow measurements: Determining the flow by using a rotameter, orifice plate and Venturi. 3) Liquid flow from a cylindrical tank: verification of a mathematical model based on the conservation equations of mass and energy to determine the time of emptying reservoir and compared the results with experimental data. 4) Determination of friction factor in pipes: evaluation of the friction factor as a function of Reynolds number in pipes. 5) Rheological tests: with the use of various types of viscometers to determine the dynamic and kinematic viscosities."

# ------------------------------------------------------------------
# Row 17 - label becomes "Avaliação:", B/C cleared entirely
# ------------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

# ------------------------------------------------------------------
# Row 18 - label becomes "Método:", B/C newly populated with the docente name
# (B18 previously had no cell at all: fix up the style afterwards because the
# worksheet's <cols> has an ambiguous overlap on column B)
# ------------------------------------------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C18").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Row 19 - label becomes "Critério:" (values already correct)
# ------------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"

# ------------------------------------------------------------------
# Row 20 - label becomes "Norma de recuperação:" (values already correct)
# ------------------------------------------------------------------
$ws.Range("A20").Value = "Norma de recuperação:"

# ------------------------------------------------------------------
# Row 21 - label becomes "Bibliografia:" (values already correct)
# ------------------------------------------------------------------
$ws.Range("A21").Value = "Bibliografia:"

# ------------------------------------------------------------------
# Row 22 - label becomes "Requisitos:", B/C cleared entirely
# ------------------------------------------------------------------
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

# ------------------------------------------------------------------
# Row 23 - A cleared, B/C newly populated with the requisito text
# (B23 previously had no cell at all: fix up the style afterwards)
# ------------------------------------------------------------------
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)`n"
$ws.Range("B14").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Row heights (customHeight) to match the new layout
# ------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# ------------------------------------------------------------------
# Old row 24 (bibliography reference text) no longer exists -- its content
# moved up into row 23, so the row itself is deleted.
# ------------------------------------------------------------------
$ws.Range("A24").EntireRow.Delete()
